# A new weekly price-report row is inserted for "Femacal de La Calera" /
# "Bruselas (repollito)" dated 2022-07-12 (Excel serial 44754), pushing the
# existing rows 38:53 down to 39:54 (dimension grows from R53 to R54).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 38 and below down by one row, carrying formatting along.
$ws.Rows("38:38").Insert()

# Populate the newly inserted row 38 with the new observation.
$ws.Range("A38").Value = 3
$ws.Range("B38").Value = "Femacal de La Calera"
$ws.Range("C38").Value = "Coquimbo"
$ws.Range("D38").Value = 44754
$ws.Range("E38").Value = 5
$ws.Range("F38").Value = 100112035
$ws.Range("G38").Value = "Bruselas (repollito)"
$ws.Range("H38").Value = "Sin especificar"
$ws.Range("I38").Value = "Primera"
$ws.Range("J38").Value = 50
$ws.Range("K38").Value = 15000
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = 15000
$ws.Range("N38").Value = "`$/malla 15 kilos"
$ws.Range("O38").Value = "Provincia de Quillota"
$ws.Range("P38").Value = 1000
$ws.Range("Q38").Value = 15
$ws.Range("R38").Value = "Hortaliza"
